$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Raw measurement updates (row 16: coremk_or0 without AOT/LW; row 18: coremk_or1) ---
# Downstream formulas (F16/F18 and the whole "TABLE FOR THESIS" block in rows 32-43)
# recalculate automatically from these.
$ws.Range("E16").Value = 774.4
$ws.Range("H16").Value = 361.3
$ws.Range("J16").Value = 288.8
$ws.Range("L16").Value = 775

$ws.Range("E18").Value = 309.8
$ws.Range("H18").Value = 194.9

# --- New "x" marker cell next to the RAW DATA header row ---
$ws.Range("C15").Value = "x"

# --- Apply the "0.0" number format to the raw value columns of the summary table ---
$ws.Range("F32:F43").NumberFormat = "0.0"
$ws.Range("H32:H43").NumberFormat = "0.0"
$ws.Range("J32:J43").NumberFormat = "0.0"
$ws.Range("L32:L43").NumberFormat = "0.0"

# --- New footer note row, highlighted in red ---
$ws.Range("A44").Value = "UPDATED 20180301"
$ws.Range("A44:B44").Interior.Color = 255

# --- Leave the cursor where the author left it ---
$ws.Range("D50").Select()
